# IHE ITI XCA Deferred - Comment Disposition
# "Final versions after TI F2F, for final line-by-line review"
#
# Fills in the Resolution ("Group") column for rows that had been left
# blank, and updates a handful of Response/Resolution dispositions on the
# "Comment Form" sheet. Also nudges the saved window/view state (zoom,
# scroll position, selection) to match the reviewer's final screen layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comment Form")

# --- Column K ("Group" / Resolution) : mark newly-reviewed rows ---------
$ws.Range("K19").Value = "Reviewed"
$ws.Range("K21").Value = "Reviewed - it is what it is; any improvement in readability will need to come from an external informative document or guide, out of the scope of this work item."
$ws.Range("K22").Value = "Reviewed"
$ws.Range("K23").Value = "Reviewed"
$ws.Range("K24").Value = "Reviewed"
$ws.Range("K25").Value = "Reviewed"
$ws.Range("K27").Value = "Reviewed"
$ws.Range("K28").Value = "Reviewed"
$ws.Range("K29").Value = "Reviewed"
$ws.Range("K30").Value = "Reviewed"
$ws.Range("K32").Value = "Reviewed"
$ws.Range("K33").Value = "Reviewed"
$ws.Range("K35").Value = "Reviewed"
$ws.Range("K36").Value = "Reviewed"
$ws.Range("K37").Value = "Drafted text in closed issues (had to make some assumptions about actual tech proposal), reviewed"

# --- Column J (Resolution status) : updated dispositions ---------------
$ws.Range("J25").Value = "No change"
$ws.Range("J27").Value = "Clarified"
$ws.Range("J33").Value = "Fixed"
$ws.Range("J37").Value = "Fixed"

# --- Window / view state, matching the reviewer's final screen ---------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 7
$win.Zoom = 150
$ws.Range("L37").Select()

$excelWin = $wb.Windows.Item(1)
$excelWin.Left = 440
$excelWin.Top = 1040
$excelWin.Width = 28000
$excelWin.Height = 14320
